# Weekly fruit/vegetable data update:
# Insert a new weekly observation as row 150 (pushing the existing rows
# 150-171 down to 151-172), then populate the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 150:171 down by one row to make room for the new record.
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new weekly record.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 45212
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112031
$ws.Cells.Item(150, 7).Value = "Poroto verde"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 30
$ws.Cells.Item(150, 11).Value = 28000
$ws.Cells.Item(150, 12).Value = 28000
$ws.Cells.Item(150, 13).Value = 28000
$ws.Cells.Item(150, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(150, 15).Value = "Perú"
$ws.Cells.Item(150, 16).Value = 1120
$ws.Cells.Item(150, 17).Value = 25
$ws.Cells.Item(150, 18).Value = "Hortaliza"
